$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Rewrite the "约瑟夫环问题描述" paragraph:
#      - bold the heading "约瑟夫环问题描述"
#      - drop the paragraph's own pPr/rPr (w:rFonts hint) mark
#      - split "...开始报数" into "...开始" + "从1顺次" + "报数..."
#      - split "...后面数到T" into "...后面" + "顺次" + "数到T..."
#      - the trailing _GoBack bookmark is dropped here (it is re-added
#        to the document's very first paragraph below)
# -----------------------------------------------------------------------
$josephXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="002602BE" w:rsidRDefault="002602BE" w:rsidP="002602BE"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>约</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>瑟</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>夫</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>环</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/></w:rPr><w:t>问题描述</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>：已知</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>num</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramStart"/><w:r><w:t>个</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>小孩（以编号1，2，3...</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>num</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>分别表示）围坐在一张圆桌周围。从编号为k的人开始</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>从1顺次</w:t></w:r><w:r><w:t>报数，数到T的那个人出列；他的下一个人又从1开始报数，后面</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>顺次</w:t></w:r><w:r><w:t>数到T的那个人又出列；依此规律重复下去，直到圆桌周围的人全部出列,试计算最后出列的那个小孩的编号。</w:t></w:r></w:p>
'@

$r = $d.Content
$found = $r.Find.Execute("约瑟夫问题描述：已知", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the Joseph-problem paragraph"
}
$r.Expand(4)
$r.InsertXML($josephXml)

# -----------------------------------------------------------------------
# 2) First embedded Visio OLE object (the flow-chart picture):
#      - shrink the preview width from 264.5pt to 264.4pt
#      - re-mint the OLEObject's ObjectID
# -----------------------------------------------------------------------
$ole1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" w:rsidR="002602BE" w:rsidRDefault="00A3799B" w:rsidP="00CA2653"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:object w:dxaOrig="6084" w:dyaOrig="4321"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:264.4pt;height:187.5pt" o:ole=""><v:imagedata r:id="rId6" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Visio.Drawing.15" ShapeID="_x0000_i1025" DrawAspect="Content" ObjectID="_1606167944" r:id="rId7"/></w:object></w:r></w:p>
'@

$r = $d.Content
$found = $r.Find.Execute("求差操作", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the anchor text before the first OLE object"
}
$r.Expand(4)
$oleRange = $d.Range($r.End, $r.End)
$oleRange.Expand(4)
$oleRange.InsertXML($ole1Xml)

# -----------------------------------------------------------------------
# 3) Second embedded Visio OLE object (the pseudo-code picture): just
#    re-mint the OLEObject's ObjectID.
# -----------------------------------------------------------------------
$ole2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" w:rsidR="00020330" w:rsidRDefault="00F73946" w:rsidP="00020330"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:object w:dxaOrig="12169" w:dyaOrig="1752"><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:445.5pt;height:64.5pt" o:ole=""><v:imagedata r:id="rId8" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Visio.Drawing.15" ShapeID="_x0000_i1026" DrawAspect="Content" ObjectID="_1606167945" r:id="rId9"/></w:object></w:r></w:p>
'@

$r = $d.Content
$found = $r.Find.Execute("问题求解", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the anchor text before the second OLE object"
}
$r.Expand(4)
$oleRange2 = $d.Range($r.End, $r.End)
$oleRange2.Expand(4)
$oleRange2 = $d.Range($oleRange2.End, $oleRange2.End)
$oleRange2.Expand(4)
$oleRange2.InsertXML($ole2Xml)

# -----------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the end of the Joseph-problem
#    paragraph (removed above in step 1) to the document's very first
#    (empty) paragraph - re-defining a bookmark with the same name moves
#    it, matching Word's own "last edit position" bookkeeping.
# -----------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1).Range
$d.Bookmarks.Add("_GoBack", $firstPara)

Write-Output "done"
